$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.861.60'
$ws.Range('E2').Value = '  +0.94%  '

$ws.Range('D3').Value = '3.585.73'
$ws.Range('E3').Value = '  +1.24%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '604.03'
$ws.Range('E5').Value = '  +1.16%  '

$ws.Range('D6').Value = '137.82'
$ws.Range('E6').Value = '  -0.70%  '

$ws.Range('D7').Value = '3.584.84'
$ws.Range('E7').Value = '  +1.20%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  +1.01%  '

$ws.Range('E10').Value = '  +0.61%  '

$ws.Range('D11').Value = '7.24'
$ws.Range('E11').Value = '  +5.02%  '

$ws.Range('D12').Value = '0.393'
$ws.Range('E12').Value = '  +1.52%  '

$ws.Range('D13').Value = '4.191.80'
$ws.Range('E13').Value = '  +1.18%  '

$ws.Range('D14').Value = '28.36'
$ws.Range('E14').Value = '  +3.86%  '

$ws.Range('D15').Value = '0.0000187'
$ws.Range('E15').Value = '  +0.58%  '

$ws.Range('D16').Value = '3.581.02'
$ws.Range('E16').Value = '  +1.01%  '

$ws.Range('D17').Value = '0.117'
$ws.Range('E17').Value = '  -0.34%  '

$ws.Range('D18').Value = '65.890.88'
$ws.Range('E18').Value = '  +0.97%  '

$ws.Range('D19').Value = '10.12'
$ws.Range('E19').Value = '  -1.37%  '

$ws.Range('E20').Value = '  +2.22%  '

$ws.Range('D21').Value = '5.88'
$ws.Range('E21').Value = '  -1.25%  '

$ws.Range('D22').Value = '395.67'
$ws.Range('E22').Value = '  +0.40%  '

$ws.Range('E23').Value = '  +2.73%  '

$ws.Range('D24').Value = '3.729.73'
$ws.Range('E24').Value = '  +1.23%  '

$ws.Range('E25').Value = '  +0.44%  '

$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').Value = '0.0000118'
$ws.Range('E27').Value = '  +2.11%  '

$ws.Range('D28').Value = '8.16'
$ws.Range('E28').Value = '  +4.36%  '

$ws.Range('E29').Value = '  +26.32%  '

$ws.Range('E30').Value = '  +2.95%  '

$ws.Range('D31').Value = '8.61'
$ws.Range('E31').Value = '  +5.69%  '

$ws.Range('E32').Value = '  -0.20%  '

$ws.Range('D33').Value = '3.589.22'
$ws.Range('E33').Value = '  +1.00%  '

$ws.Range('D34').Value = '24.55'
$ws.Range('E34').Value = '  +2.89%  '

$ws.Range('E35').Value = '  +1.87%  '

$ws.Range('E36').Value = '  -0.01%  '

$ws.Range('D37').Value = '5.40'
$ws.Range('E37').Value = '  +7.82%  '

$ws.Range('D38').Value = '1.64'
$ws.Range('E38').Value = '  +5.27%  '

$ws.Range('D39').Value = '7.05'
$ws.Range('E39').Value = '  +1.07%  '

$ws.Range('D40').Value = '168.07'
$ws.Range('E40').Value = '  -0.69%  '

$ws.Range('D41').Value = '0.0839'
$ws.Range('E41').Value = '  +4.45%  '

$ws.Range('E42').Value = '  +1.36%  '

$ws.Range('D43').Value = '27.00'
$ws.Range('E43').Value = '  +2.37%  '

$ws.Range('E44').Value = '  +7.35%  '

$ws.Range('D45').Value = '43.14'
$ws.Range('E45').Value = '  +0.87%  '

$ws.Range('D46').Value = '4.55'
$ws.Range('E46').Value = '  +2.77%  '

$ws.Range('D47').Value = '0.999'
$ws.Range('E47').Value = '  -0.10%  '

$ws.Range('D48').Value = '1.71'
$ws.Range('E48').Value = '  +1.85%  '

$ws.Range('D49').Value = '7.02'
$ws.Range('E49').Value = '  +3.07%  '

$ws.Range('D50').Value = '2.464.20'
$ws.Range('E50').Value = '  +2.80%  '

$ws.Range('E51').Value = '  +4.74%  '
